# Insert a new weekly price record as row 324 on the single worksheet,
# pushing the existing rows 324:369 down to 325:370 (dimension becomes A1:R370).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 324:369 down by one row, carrying formatting (same as
# right-click > Insert on a whole row in Excel).
$ws.Rows.Item(324).Insert()

# Populate the newly inserted row 324 with the new record.
$ws.Range("A324").Value = 9
$ws.Range("B324").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C324").Value = "Metropolitana"
$ws.Range("D324").Value = 44637
$ws.Range("E324").Value = 13
$ws.Range("F324").Value = 100112013
$ws.Range("G324").Value = "Alcachofa"
$ws.Range("H324").Value = "Española"
$ws.Range("I324").Value = "Primera"
$ws.Range("J324").Value = 50
$ws.Range("K324").Value = 20000
$ws.Range("L324").Value = 20000
$ws.Range("M324").Value = 20000
$ws.Range("N324").Value = "$/caja 30 unidades"
$ws.Range("O324").Value = "Provincia de Limarí"
$ws.Range("P324").Value = 667
$ws.Range("Q324").Value = 30
$ws.Range("R324").Value = "Hortaliza"
